{"js": "// Replace the division-problem text in each table cell with its new\n// value. Each old value is unique in the document, so a matchCase,\n// matchWholeWord search on the body finds exactly one hit per pair and\n// insertText(..., Word.InsertLocation.replace) swaps the run's text in\n// place, preserving the existing run formatting (font/size).\nconst replacements = [\n  [\"95\u00f74=\", \"77\u00f72=\"],\n  [\"89\u00f76=\", \"91\u00f78=\"],\n  [\"74\u00f73=\", \"65\u00f78=\"],\n  [\"28\u00f74=\", \"51\u00f73=\"],\n  [\"23\u00f75=\", \"78\u00f73=\"],\n  [\"42\u00f76=\", \"36\u00f79=\"],\n  [\"20\u00f77=\", \"87\u00f73=\"],\n  [\"93\u00f77=\", \"58\u00f75=\"],\n  [\"64\u00f77=\", \"39\u00f75=\"],\n  [\"55\u00f75=\", \"56\u00f72=\"],\n  [\"26\u00f76=\", \"59\u00f73=\"],\n  [\"64\u00f74=\", \"93\u00f77=\"],\n  [\"11\u00f78=\", \"89\u00f78=\"],\n  [\"90\u00f72=\", \"23\u00f74=\"],\n  [\"90\u00f76=\", \"72\u00f79=\"],\n  [\"62\u00f77=\", \"86\u00f74=\"],\n  [\"99\u00f78=\", \"50\u00f72=\"],\n  [\"24\u00f79=\", \"80\u00f79=\"],\n  [\"61\u00f78=\", \"77\u00f72=\"],\n  [\"68\u00f75=\", \"48\u00f74=\"],\n  [\"14\u00f78=\", \"70\u00f79=\"],\n  [\"41\u00f78=\", \"68\u00f73=\"],\n  [\"33\u00f78=\", \"36\u00f76=\"],\n  [\"16\u00f79=\", \"34\u00f76=\"],\n  [\"92\u00f74=\", \"21\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division-problem text in each table cell with its new\n# value. Each old value is unique in the document, so Find/Replace\n# against a fresh whole-document range per pair finds exactly one hit\n# and swaps the text while leaving run formatting (font/size) intact.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"95\u00f74=\", \"77\u00f72=\"),\n  @(\"89\u00f76=\", \"91\u00f78=\"),\n  @(\"74\u00f73=\", \"65\u00f78=\"),\n  @(\"28\u00f74=\", \"51\u00f73=\"),\n  @(\"23\u00f75=\", \"78\u00f73=\"),\n  @(\"42\u00f76=\", \"36\u00f79=\"),\n  @(\"20\u00f77=\", \"87\u00f73=\"),\n  @(\"93\u00f77=\", \"58\u00f75=\"),\n  @(\"64\u00f77=\", \"39\u00f75=\"),\n  @(\"55\u00f75=\", \"56\u00f72=\"),\n  @(\"26\u00f76=\", \"59\u00f73=\"),\n  @(\"64\u00f74=\", \"93\u00f77=\"),\n  @(\"11\u00f78=\", \"89\u00f78=\"),\n  @(\"90\u00f72=\", \"23\u00f74=\"),\n  @(\"90\u00f76=\", \"72\u00f79=\"),\n  @(\"62\u00f77=\", \"86\u00f74=\"),\n  @(\"99\u00f78=\", \"50\u00f72=\"),\n  @(\"24\u00f79=\", \"80\u00f79=\"),\n  @(\"61\u00f78=\", \"77\u00f72=\"),\n  @(\"68\u00f75=\", \"48\u00f74=\"),\n  @(\"14\u00f78=\", \"70\u00f79=\"),\n  @(\"41\u00f78=\", \"68\u00f73=\"),\n  @(\"33\u00f78=\", \"36\u00f76=\"),\n  @(\"16\u00f79=\", \"34\u00f76=\"),\n  @(\"92\u00f74=\", \"21\u00f73=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $rng = $d.Content\n  $find = $rng.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
